$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New status "Partial" added next to "Makefile" row (G4), referencing the
# new shared string "Partial".
$ws.Range("G4").Value = "Partial"

# Error Handling, Logging, Docker deploy + Compose, Makefile now checked.
$ws.Range("B11").Value = $true
$ws.Range("B12").Value = $true
$ws.Range("B13").Value = $true
$ws.Range("B14").Value = $true

# Unit Tests row gains an "Optional" note in column C.
$ws.Range("C22").Value = "Optional"

# Update the active selection to reflect where the author was working.
$ws.Range("E17").Select()
